$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add three new PPC datapoint abbreviation rows to the device datapoints
# table:
#   - PPC_P_SET_MODE          inserted before PPC_P_SET_REL
#   - PPC_Q_SET_MODE          inserted before PPC_Q_SET_REL
#   - PPC_V_REF_Q_V_SHIFT     inserted after  PPC_Q_SET_REL
# ---------------------------------------------------------------------------

$phi = [char]0x03C6

$pSetModeDesc = "Active power control method  0: No configuration found 1: Variable fixed value Pvar fix 2: Variable fixed value Pvar DI 3: Variable fixed value Pvar AI 4: Variable fixed value Pvar Modbus 5: Remote Power Control (RPC) 100: LFSM-O 101: LFSM-U 102: FSM 112: RPC (Remote Power Control) & FSM 200: Fail-safe operation (hold last setpoint) 201: Fail-safe operation (default setpoint) 202: Fail-safe operation (system fallback value) 203: Fail-safe operation (Automatic grid disconnection)"

$qSetModeDesc = "Reactive power control method  0: No configuration found 1: Variable fixed value cos " + $phi + "var fix 2: Variable fixed value cos " + $phi + "var DI 3: Variable fixed value cos " + $phi + "var AI 4: Variable fixed value cos " + $phi + "var Modbus 5: Characteristic curve cos " + $phi + " (P) 6: Characteristic curve cos " + $phi + " (V) 7: Variable fixed value Qvar fix 8: Variable fixed value Qvar DI 9: Variable fixed value Qvar AI 10: Variable fixed value Qvar Modbus 11: Characteristic curve Q (P) 12: Characteristic curve Q (V) 13: Characteristic curve Q (tan " + $phi + ") 14: Voltage control Q (V droop) 15: Characteristic curve Q(V) Modbus, from firmware 28.0.2 on 16: Characteristic curve cos " + $phi + " (P) Modbus, from firmware 28.0.2 on 100: Reactive power compensation 200: Fail-safe operation (hold last setpoint) 201: Fail-safe operation (default setpoint) 202: Fail-safe operation (system fallback value)"

# Row 35: PPC_P_SET_MODE (inserted above PPC_P_SET_REL)
$ws.Rows(35).Insert()
$ws.Range("A35").Value = "datapoints"
$ws.Range("B35").Value = "PPC_P_SET_MODE"
$ws.Range("D35").Value = $pSetModeDesc

# Row 44: PPC_Q_SET_MODE (inserted above PPC_Q_SET_REL)
$ws.Rows(44).Insert()
$ws.Range("A44").Value = "datapoints"
$ws.Range("B44").Value = "PPC_Q_SET_MODE"
$ws.Range("D44").Value = $qSetModeDesc

# Row 46: PPC_V_REF_Q_V_SHIFT (inserted below PPC_Q_SET_REL)
$ws.Rows(46).Insert()
$ws.Range("A46").Value = "datapoints"
$ws.Range("B46").Value = "PPC_V_REF_Q_V_SHIFT"
$ws.Range("C46").Value = "V"
$ws.Range("D46").Value = "Voltage shift for Q(V) curve"

# The Description column grows very wide to fit the new control-method text.
$ws.Columns("D").ColumnWidth = 987.64
